$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the cell to be stored as literal text so that values which
    # look numeric (e.g. "7.10", "300.52") keep their exact original
    # formatting instead of being auto-converted to a number by Excel.
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.ClearFormats()
}

$ws.Range("D2").Value = "45.919.17"
$ws.Range("E2").Value = "  -2.08%  "
$ws.Range("D3").Value = "2.379.77"
$ws.Range("E3").Value = "  +3.14%  "
$ws.Range("E4").Value = "  +0.02%  "
Set-TextValue "D5" "300.52"
$ws.Range("E5").Value = "  -1.13%  "
Set-TextValue "D6" "98.36"
$ws.Range("E6").Value = "  -3.41%  "
$ws.Range("E7").Value = "  -0.66%  "
$ws.Range("E8").Value = "  +0.08%  "
Set-TextValue "D9" "0.507"
$ws.Range("E9").Value = "  -4.52%  "
Set-TextValue "D10" "34.44"
$ws.Range("E10").Value = "  -6.76%  "
$ws.Range("E11").Value = "  -2.13%  "
Set-TextValue "D12" "7.10"
$ws.Range("E12").Value = "  -5.16%  "
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("D14").Value = "2.739.88"
$ws.Range("E14").Value = "  +3.07%  "
$ws.Range("D15").Value = "2.399.57"
$ws.Range("E15").Value = "  +4.02%  "
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("E17").Value = "  -2.09%  "
$ws.Range("D18").Value = "45.880.30"
$ws.Range("E18").Value = "  -2.07%  "
$ws.Range("E19").Value = "  -7.55%  "
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("E21").Value = "  -1.44%  "
Set-TextValue "D22" "66.64"
$ws.Range("E22").Value = "  -0.51%  "
Set-TextValue "D23" "242.74"
$ws.Range("E23").Value = "  -2.62%  "
$ws.Range("E24").Value = "  -5.57%  "
Set-TextValue "D25" "0.999"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").Value = "  -2.39%  "
Set-TextValue "D27" "39.04"
$ws.Range("E27").Value = "  -11.44%  "
Set-TextValue "D28" "2.21"
$ws.Range("E28").Value = "  -3.14%  "
Set-TextValue "D29" "9.71"
$ws.Range("E29").Value = "  -2.46%  "
Set-TextValue "D30" "20.91"
$ws.Range("E30").Value = "  +3.25%  "
$ws.Range("E31").Value = "  +16.91%  "
$ws.Range("E32").Value = "  +6.66%  "
$ws.Range("E33").Value = "  -4.85%  "
Set-TextValue "D34" "146.16"
$ws.Range("E34").Value = "  -0.81%  "
Set-TextValue "D35" "0.0770"
$ws.Range("E35").Value = "  -3.97%  "
Set-TextValue "D36" "0.112"
$ws.Range("E36").Value = "  -1.91%  "
Set-TextValue "D37" "1.91"
$ws.Range("E37").Value = "  +5.60%  "
$ws.Range("E38").Value = "  -2.90%  "
Set-TextValue "D39" "14.86"
$ws.Range("E39").Value = "  -7.83%  "
Set-TextValue "D40" "3.86"
$ws.Range("E40").Value = "  -3.92%  "
Set-TextValue "D41" "0.0298"
$ws.Range("E41").Value = "  -3.07%  "
Set-TextValue "D42" "3.19"
$ws.Range("E42").Value = "  -8.66%  "
$ws.Range("D43").Value = "1.937.11"
$ws.Range("E43").Value = "  +3.78%  "
$ws.Range("E44").Value = "  +0.13%  "
Set-TextValue "D45" "93.40"
$ws.Range("E45").Value = "  +5.42%  "
$ws.Range("E46").Value = "  -10.11%  "
Set-TextValue "D47" "8.48"
$ws.Range("E47").Value = "  +4.97%  "
$ws.Range("E48").Value = "  -6.38%  "
Set-TextValue "D49" "99.12"
$ws.Range("E49").Value = "  +1.85%  "
$ws.Range("D50").Value = "2.612.66"
$ws.Range("E50").Value = "  +3.07%  "
Set-TextValue "D51" "68.85"
$ws.Range("E51").Value = "  -8.19%  "
